# Update cryptos list price/volume data to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.280.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.281.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.89"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.949"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.629.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.278.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.078.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +29.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.33"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.63"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0865"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -12.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.75"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.23%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.06"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.605.72"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.26%  "
